# SubRES_REZoning_Sol_Win_and_Hydro.xlsx edit
#
# The "solar" sheet has a wind-grid-connection table in columns R:Y / AC:AG
# (rows 4-27), one row per CHE_n grid cell (process "distr_elc_won-CHE_000N"
# connecting solar+wind resources to buses in that cell). The underlying source
# data was regenerated upstream with a new row assignment, so every row's
# process/description/comm-in/comm-out/grid_cell got reassigned to a different
# CHE_n cell (values for columns U,V,W,X,Y,AF stay constant across all rows, so
# only S,T,AC,AD,AE,AG need to change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

# Row => @(process, description, comm-in, comm-out, grid_cell)
$rowData = @{
    4 = @('distr_elc_won-CHE_0000', 'connecting solar and wind to buses in grid cell CHE_0', 'elc_won-CHE_0000,elc_spv-CHE_0000', 'e_CH51-220,e_CH52-220,e_CH57-220,e_w177392130-220,e_w239937062-220', 'CHE_0')
    5 = @('distr_elc_won-CHE_0020', 'connecting solar and wind to buses in grid cell CHE_20', 'elc_won-CHE_0020,elc_spv-CHE_0020', 'e_CH21-220,e_CH22-220,e_CH29-220,e_CH41-380,e_w1208713169-220,e_w207993342-220,e_w208780268-380,e_w212498548-220,e_w36348118-220,e_w365556107-220,e_w71500123-220', 'CHE_20')
    6 = @('distr_elc_won-CHE_0001', 'connecting solar and wind to buses in grid cell CHE_1', 'elc_won-CHE_0001,elc_spv-CHE_0001', 'e_CH31-220,e_w132373704-220,e_w232662311-220,e_w44496892-220,e_w55695765-220', 'CHE_1')
    7 = @('distr_elc_won-CHE_0006', 'connecting solar and wind to buses in grid cell CHE_6', 'elc_won-CHE_0006,elc_spv-CHE_0006', 'e_w127004407-380', 'CHE_6')
    8 = @('distr_elc_won-CHE_0011', 'connecting solar and wind to buses in grid cell CHE_11', 'elc_won-CHE_0011,elc_spv-CHE_0011', 'e_CH46-220,e_CH47-220,e_w228003081-220,e_w391576135-220,e_w391577741-220,e_w969819301-220', 'CHE_11')
    9 = @('distr_elc_won-CHE_0015', 'connecting solar and wind to buses in grid cell CHE_15', 'elc_won-CHE_0015,elc_spv-CHE_0015', 'e_CH50-220,e_CH56-220,e_CH58-220,e_CH59-220,e_w1327084723-220,e_w281800404-220,e_w281803398-220,e_w281815404-220,e_w35487135-220', 'CHE_15')
    10 = @('distr_elc_won-CHE_0025', 'connecting solar and wind to buses in grid cell CHE_25', 'elc_won-CHE_0025,elc_spv-CHE_0025', 'e_r7933294-380', 'CHE_25')
    11 = @('distr_elc_won-CHE_0010', 'connecting solar and wind to buses in grid cell CHE_10', 'elc_won-CHE_0010,elc_spv-CHE_0010', 'e_CH11-220,e_w109037817-220,e_w127004407-380,e_w27435934-220,e_w30350721-220,e_w397960460-380,e_w88901626-380', 'CHE_10')
    12 = @('distr_elc_won-CHE_0022', 'connecting solar and wind to buses in grid cell CHE_22', 'elc_won-CHE_0022,elc_spv-CHE_0022', 'e_CH4-220,e_CH9-220,e_w455120191-220,e_w83861269-220,e_w92798668-220,e_w92873516-220', 'CHE_22')
    13 = @('distr_elc_won-CHE_0024', 'connecting solar and wind to buses in grid cell CHE_24', 'elc_won-CHE_0024,elc_spv-CHE_0024', 'e_CH30-380,e_CH33-380,e_CH37-380,e_r7933294-380', 'CHE_24')
    14 = @('distr_elc_won-CHE_0008', 'connecting solar and wind to buses in grid cell CHE_8', 'elc_won-CHE_0008,elc_spv-CHE_0008', 'e_CH43-220,e_w209324991-220,e_w26166640-220,e_w402055336-220,e_w758315582-220', 'CHE_8')
    15 = @('distr_elc_won-CHE_0005', 'connecting solar and wind to buses in grid cell CHE_5', 'elc_won-CHE_0005,elc_spv-CHE_0005', 'e_r5378910-220,e_w161853746-220', 'CHE_5')
    16 = @('distr_elc_won-CHE_0014', 'connecting solar and wind to buses in grid cell CHE_14', 'elc_won-CHE_0014,elc_spv-CHE_0014', 'e_CH1-220,e_CH2-220,e_CH5-220,e_CH7-220,e_w108257952-220,e_w1284913429-220,e_w190819048-220,e_w27107779-220,e_w31308888-220,e_w33271433-220,e_w356292116-220,e_w35840165-380,e_w50319857-220,e_w89405664-220,e_w89977424-220,e_w98648381-220', 'CHE_14')
    17 = @('distr_elc_won-CHE_0018', 'connecting solar and wind to buses in grid cell CHE_18', 'elc_won-CHE_0018,elc_spv-CHE_0018', 'e_CH3-220,e_CH6-220,e_r9310861-220,e_w11282314-220,e_w147714395-220,e_w148015471-220,e_w22899676-220,e_w240575085-220,e_w26843160-220,e_w281809991-220,e_w50561341-220,e_w87281514-220', 'CHE_18')
    18 = @('distr_elc_won-CHE_0003', 'connecting solar and wind to buses in grid cell CHE_3', 'elc_won-CHE_0003,elc_spv-CHE_0003', 'e_CH48-225,e_CH49-225,e_CH53-225,e_CH60-225,e_w234983117-220,e_w238138373-380,e_w260211728-225,e_w55698557-220,e_w802058337-220', 'CHE_3')
    19 = @('distr_elc_won-CHE_0007', 'connecting solar and wind to buses in grid cell CHE_7', 'elc_won-CHE_0007,elc_spv-CHE_0007', 'e_CH44-220,e_w212722603-220,e_w236819191-220,e_w758943072-220', 'CHE_7')
    20 = @('distr_elc_won-CHE_0017', 'connecting solar and wind to buses in grid cell CHE_17', 'elc_won-CHE_0017,elc_spv-CHE_0017', 'e_CH12-220,e_CH13-220,e_CH16-380,e_CH18-220,e_w192677427-220,e_w52738225-220', 'CHE_17')
    21 = @('distr_elc_won-CHE_0019', 'connecting solar and wind to buses in grid cell CHE_19', 'elc_won-CHE_0019,elc_spv-CHE_0019', 'e_CH45-220,e_w281804158-220', 'CHE_19')
    22 = @('distr_elc_won-CHE_0023', 'connecting solar and wind to buses in grid cell CHE_23', 'elc_won-CHE_0023,elc_spv-CHE_0023', 'e_w207991759-380', 'CHE_23')
    23 = @('distr_elc_won-CHE_0013', 'connecting solar and wind to buses in grid cell CHE_13', 'elc_won-CHE_0013,elc_spv-CHE_0013', 'e_CH15-220,e_w146225999-220,e_w159527493-220,e_w242269161-220,e_w281799252-220,e_w35002638-220,e_w97941869-220', 'CHE_13')
    24 = @('distr_elc_won-CHE_0009', 'connecting solar and wind to buses in grid cell CHE_9', 'elc_won-CHE_0009,elc_spv-CHE_0009', 'e_CH14-220,e_w1105061707-220,e_w147557680-220,e_w165254212-220,e_w402053379-220', 'CHE_9')
    25 = @('distr_elc_won-CHE_0021', 'connecting solar and wind to buses in grid cell CHE_21', 'elc_won-CHE_0021,elc_spv-CHE_0021', 'e_CH17-380,e_w211907009-220', 'CHE_21')
    26 = @('distr_elc_won-CHE_0004', 'connecting solar and wind to buses in grid cell CHE_4', 'elc_won-CHE_0004,elc_spv-CHE_0004', 'e_w111162936-220,e_w122720993-220,e_w194258388-220', 'CHE_4')
    27 = @('distr_elc_won-CHE_0012', 'connecting solar and wind to buses in grid cell CHE_12', 'elc_won-CHE_0012,elc_spv-CHE_0012', 'e_CH19-220,e_CH20-220,e_CH25-220,e_CH27-220,e_CH34-220,e_CH35-220,e_CH36-220,e_CH38-220,e_CH39-220,e_CH40-220,e_w1086214433-220,e_w1092884227-220,e_w240959264-220,e_w364949845-220', 'CHE_12')
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("S$row").Value = $vals[0]
    $ws.Range("T$row").Value = $vals[1]
    $ws.Range("AC$row").Value = $vals[0]
    $ws.Range("AD$row").Value = $vals[2]
    $ws.Range("AE$row").Value = $vals[3]
    $ws.Range("AG$row").Value = $vals[4]
}

